$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
